# Weekly CompStat data refresh for the 33rd Precinct report.
# Advances the reporting week from 1/1/2024-1/7/2024 (Vol 31 No 1) to
# 1/8/2024-1/14/2024 (Vol 31 No 2), and writes the newly collected crime
# counts / period-over-period percentage changes onto the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: week-ending label + volume/number ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"


# Row 16 (Robbery)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -25
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("K16").Value = 0
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = 100
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = -20

# Row 17 (Fel. Assault)
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 26.315789473684
$ws.Range("I17").Value = 14
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 27.272727272727
$ws.Range("L17").Value = 40
$ws.Range("M17").Value = 366.666666666667

# Row 18 (Burglary)
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = -28.571428571428
$ws.Range("L18").Value = -44.444444444444
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = -37.5

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 10
$ws.Range("K19").Value = -10
$ws.Range("L19").Value = 28.571428571428

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -30
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 6
$ws.Range("L20").Value = -14.285714285714

# Row 21 (TOTAL)
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = 6.060606060606
$ws.Range("I21").Value = 38
$ws.Range("J21").Value = 38
$ws.Range("L21").Value = 8.571428571428
$ws.Range("M21").Value = 90

# Row 22 (Transit)
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 2
$ws.Range("L22").Value = 100
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 23 (Housing)
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 2
$ws.Range("L23").Value = 0
$ws.Range("L23").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -20.27027027027
$ws.Range("I24").Value = 25
$ws.Range("J24").Value = 27
$ws.Range("K24").Value = -7.407407407407
$ws.Range("L24").Value = -60.9375
$ws.Range("M24").Value = 38.888888888888

# Row 25 (Misd. Assault)
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -25.714285714285
$ws.Range("I25").Value = 13
$ws.Range("J25").Value = 15
$ws.Range("K25").Value = -13.333333333333
$ws.Range("L25").Value = 8.333333333333
$ws.Range("M25").Value = -31.578947368421

# Row 26 (UCR Rape*)
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = 0
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J26").Value = 1
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = -100
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 1
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("L27").Value = 0

# Row 36 (Historical Perspective: Murder)
$ws.Range("J36").Value = 8
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 33.333333333333

# Row 38 (Historical Perspective: Robbery)
$ws.Range("J38").Value = 177
$ws.Range("K38").Value = -35.164835164835
$ws.Range("L38").Value = -45.871559633027
$ws.Range("M38").Value = 24.647887323943

# Row 42 (Historical Perspective: G.L.A.)
$ws.Range("J42").Value = 142
$ws.Range("K42").Value = -26.80412371134
$ws.Range("L42").Value = -49.82332155477
$ws.Range("M42").Value = 6.766917293233

# Row 43 (Historical Perspective: TOTAL)
$ws.Range("J43").Value = 1052
$ws.Range("K43").Value = -12.841756420878
$ws.Range("L43").Value = -36.816816816816
$ws.Range("M43").Value = 43.91244870041
